$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = [double]"0.007064987288959285"
$ws.Range("E2").Value = [double]"0.007064987288959285"

# Row 3
$ws.Range("D3").Value = [double]"3.550725809296352E-07"
$ws.Range("E3").Value = [double]"3.550725809296352E-07"

# Row 4
$ws.Range("D4").Value = [double]"1.884699157459474E-10"
$ws.Range("E4").Value = [double]"1.884699157459474E-10"

# Row 5
$ws.Range("D5").Value = [double]"0.01079066539298056"
$ws.Range("E5").Value = [double]"0.01079066539298056"

# Row 6
$ws.Range("D6").Value = [double]"0.8689145827006219"
$ws.Range("E6").Value = [double]"0.8689145827006219"

# Row 7
$ws.Range("D7").Value = [double]"0.9989297622280975"
$ws.Range("E7").Value = [double]"0.001070237771902516"

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"0.0001725185849984225"
$ws.Range("E8").Value = [double]"0.9998274814150016"

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = [double]"0.0001416071617547044"
$ws.Range("E9").Value = [double]"0.9998583928382453"

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"0.01716371028495326"
$ws.Range("E10").Value = [double]"0.9828362897150468"

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = [double]"6.313642712226358E-06"
$ws.Range("E11").Value = [double]"0.9999936863572878"
$ws.Range("F11").Value = [double]"3.561613082885742"
$ws.Range("G11").Value = [double]"0.5"

# Row 12
$ws.Range("D12").Value = [double]"0.0002018871026234869"
$ws.Range("E12").Value = [double]"0.0002018871026234869"

# Row 13
$ws.Range("D13").Value = [double]"1.655497249263877E-11"
$ws.Range("E13").Value = [double]"1.655497249263877E-11"

# Row 14
$ws.Range("D14").Value = [double]"3.318938310843055E-16"
$ws.Range("E14").Value = [double]"3.318938310843055E-16"

# Row 15
$ws.Range("D15").Value = [double]"7.440602584897592E-05"
$ws.Range("E15").Value = [double]"7.440602584897592E-05"

# Row 16
$ws.Range("D16").Value = [double]"0.908369655802864"
$ws.Range("E16").Value = [double]"0.908369655802864"

# Row 17
$ws.Range("D17").Value = [double]"0.9999864142652971"
$ws.Range("E17").Value = [double]"1.358573470289226E-05"

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"2.339118916565548E-07"
$ws.Range("E18").Value = [double]"0.9999997660881084"

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"2.469244484908635E-07"
$ws.Range("E19").Value = [double]"0.9999997530755516"

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"0.000297710976382215"
$ws.Range("E20").Value = [double]"0.9997022890236178"

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = [double]"4.614538534118149E-09"
$ws.Range("E21").Value = [double]"0.9999999953854615"
$ws.Range("F21").Value = [double]"6.018622875213623"
$ws.Range("G21").Value = [double]"0.5"

